{"js": "// Replace the date line and every \"a\u00f7b=c, d\" answer cell with its updated\n// value. Every search string below is unique in the document, so a plain\n// body.search() + replace on the (single) hit is safe and avoids any\n// cross-matching between the 26 edits.\nconst replacements = [\n  [\"2024-12-18 Wednesday\", \"2024-12-19 Thursday\"],\n  [\"115\u00f77=16, 3\", \"156\u00f76=26, 0\"],\n  [\"931\u00f74=232, 3\", \"631\u00f72=315, 1\"],\n  [\"613\u00f75=122, 3\", \"923\u00f77=131, 6\"],\n  [\"580\u00f76=96, 4\", \"659\u00f74=164, 3\"],\n  [\"230\u00f76=38, 2\", \"405\u00f72=202, 1\"],\n  [\"628\u00f79=69, 7\", \"645\u00f72=322, 1\"],\n  [\"886\u00f75=177, 1\", \"614\u00f75=122, 4\"],\n  [\"670\u00f76=111, 4\", \"452\u00f73=150, 2\"],\n  [\"373\u00f77=53, 2\", \"292\u00f76=48, 4\"],\n  [\"260\u00f73=86, 2\", \"693\u00f78=86, 5\"],\n  [\"633\u00f78=79, 1\", \"253\u00f78=31, 5\"],\n  [\"855\u00f76=142, 3\", \"423\u00f74=105, 3\"],\n  [\"118\u00f78=14, 6\", \"211\u00f73=70, 1\"],\n  [\"562\u00f77=80, 2\", \"821\u00f79=91, 2\"],\n  [\"211\u00f72=105, 1\", \"911\u00f75=182, 1\"],\n  [\"960\u00f79=106, 6\", \"428\u00f72=214, 0\"],\n  [\"445\u00f75=89, 0\", \"387\u00f72=193, 1\"],\n  [\"393\u00f74=98, 1\", \"885\u00f75=177, 0\"],\n  [\"166\u00f75=33, 1\", \"836\u00f72=418, 0\"],\n  [\"898\u00f79=99, 7\", \"443\u00f78=55, 3\"],\n  [\"816\u00f78=102, 0\", \"800\u00f75=160, 0\"],\n  [\"196\u00f77=28, 0\", \"985\u00f75=197, 0\"],\n  [\"898\u00f77=128, 2\", \"699\u00f77=99, 6\"],\n  [\"985\u00f77=140, 5\", \"332\u00f78=41, 4\"],\n  [\"888\u00f74=222, 0\", \"196\u00f72=98, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"a\u00f7b=c, d\" answer cell with its updated\n# value. Each search string occurs exactly once in the document, so a\n# single Find/Replace pass per pair (re-run over the whole document) is\n# enough and keeps each edit independent of the others.\n$replacements = @(\n    @(\"2024-12-18 Wednesday\", \"2024-12-19 Thursday\"),\n    @(\"115\u00f77=16, 3\", \"156\u00f76=26, 0\"),\n    @(\"931\u00f74=232, 3\", \"631\u00f72=315, 1\"),\n    @(\"613\u00f75=122, 3\", \"923\u00f77=131, 6\"),\n    @(\"580\u00f76=96, 4\", \"659\u00f74=164, 3\"),\n    @(\"230\u00f76=38, 2\", \"405\u00f72=202, 1\"),\n    @(\"628\u00f79=69, 7\", \"645\u00f72=322, 1\"),\n    @(\"886\u00f75=177, 1\", \"614\u00f75=122, 4\"),\n    @(\"670\u00f76=111, 4\", \"452\u00f73=150, 2\"),\n    @(\"373\u00f77=53, 2\", \"292\u00f76=48, 4\"),\n    @(\"260\u00f73=86, 2\", \"693\u00f78=86, 5\"),\n    @(\"633\u00f78=79, 1\", \"253\u00f78=31, 5\"),\n    @(\"855\u00f76=142, 3\", \"423\u00f74=105, 3\"),\n    @(\"118\u00f78=14, 6\", \"211\u00f73=70, 1\"),\n    @(\"562\u00f77=80, 2\", \"821\u00f79=91, 2\"),\n    @(\"211\u00f72=105, 1\", \"911\u00f75=182, 1\"),\n    @(\"960\u00f79=106, 6\", \"428\u00f72=214, 0\"),\n    @(\"445\u00f75=89, 0\", \"387\u00f72=193, 1\"),\n    @(\"393\u00f74=98, 1\", \"885\u00f75=177, 0\"),\n    @(\"166\u00f75=33, 1\", \"836\u00f72=418, 0\"),\n    @(\"898\u00f79=99, 7\", \"443\u00f78=55, 3\"),\n    @(\"816\u00f78=102, 0\", \"800\u00f75=160, 0\"),\n    @(\"196\u00f77=28, 0\", \"985\u00f75=197, 0\"),\n    @(\"898\u00f77=128, 2\", \"699\u00f77=99, 6\"),\n    @(\"985\u00f77=140, 5\", \"332\u00f78=41, 4\"),\n    @(\"888\u00f74=222, 0\", \"196\u00f72=98, 0\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Search text not found: $oldText\"\n    }\n}\n"}
